$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cell = $ws.Range("L6")
$cell.Value = $cell.Value2 + "|crystalball;2"
